$d = $word.ActiveDocument
Write-Output ("Paragraphs: " + $d.Paragraphs.Count)
Write-Output ("Tables: " + $d.Tables.Count)
$t = $d.Tables.Item(1)
Write-Output ("Rows: " + $t.Rows.Count)
Write-Output ("Cols: " + $t.Columns.Count)
$cell = $t.Cell(1,3)
$r = $cell.Range
Write-Output "---CELL XML---"
Write-Output $r.WordOpenXML
